# Applies the Case_10_1 pl_mw.xlsx update:
# - Extends the table from columns B:O to B:Q (adds two new columns, headers 14 and 15)
# - Recomputes all the data rows (2-25) for columns B-I and O, and fills the new P/Q columns
# - Columns J:N remain all zeros (unchanged), column A (index) and the row-1 header values
#   for B:O are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend the header row (row 1) with the two new columns P and Q ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy the header formatting (bold font, thin border, centered) from O1 onto the
# two new header cells so they look like the rest of the header row.
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- 2. Rewrite the data rows (row 2 .. row 25) ---
# Columns, in the order the embedded data appears below, are: B C D E F G H I O P Q
$columns = @("B","C","D","E","F","G","H","I","O","P","Q")

$csv = @"
3.46677384630982,1.019169272890139,0.04824190763466873,1.35169962688731,0.5084535737087279,0.0007896385517037115,0.01092094002687105,0.003536527785124033,0,0,1.435151477765999
3.02011174815118,0.8990605546520669,0.0433835340809452,1.176740016824183,0.468187038062311,0.000793366441376707,0.00766975902495437,0.001965288275290966,0,0,1.348594257482659
2.745669458445775,0.8258832939075091,0.04041041276209967,1.069666178277558,0.4441198563540496,0.0007957272508090085,0.005904022270549425,0.001270814676121557,0,0,1.297464469391542
2.633730900298815,0.7977722387358028,0.03928077046956702,1.02608314129057,0.4336700664951891,0.0007967124296093884,0.00523637476029204,0.001110499418202515,0,0,1.274436995028537
2.615103766652396,0.7950489896263662,0.03919122151063448,1.018830727099584,0.4309748955527013,0.0007968828780179578,0.005126972417599962,0.001168165231651308,0,0,1.267397047251904
2.74406360535454,0.8308092758339853,0.04066341026042863,1.069026406947373,0.4413300423773592,0.0007957557770997617,0.005890106154369357,0.001472853056091239,0,0,1.28828808501018
3.312653545137607,0.9847298999007421,0.04692401796475565,1.291216389130142,0.4909136236972671,0.0007909279253908628,0.009741767389663553,0.003174773938467546,0,0,1.393125671300481
4.427610317638539,1.282110248089054,0.0588090931990024,1.730915313986756,0.5988199823157387,0.0007819898167375172,0.01935802986678037,0.008607690668569923,0,0,1.632718893241275
5.24749449069725,1.50644471872306,0.06925275944912102,1.955498156088424,0.6683724546327312,0.0007759177289413412,0.02725836382229918,0.0143664161419883,0,0,1.775547518866148
5.612142301903873,1.610434802784823,0.08861588644301577,1.267337406619561,0.5871507068896307,0.0007746196692268896,0.04270868006982909,0.01643561128637661,0,0,1.474547141099578
5.74725555498668,1.642389293035194,0.1037805964393641,0.7713189373199754,0.5121260620343264,0.0007745905169331831,0.07854308646568597,0.01664005842228544,0,0,1.224683832186997
5.71206086016997,1.628934670659646,0.1168170525753141,0.3948377359290731,0.4334589186537912,0.0007755528712741501,0.1314041879758321,0.01567717074598907,0,0,0.9858811040282944
5.610402796626545,1.600782500645209,0.1251548963565625,0.2040035724406977,0.3772055239524477,0.0007766576504858142,0.1788295655706946,0.01459786682350295,0,0,0.8260307867355579
5.549222564864976,1.586479300100507,0.1266680095887693,0.1664673984577796,0.361549071281118,0.0007771664842525621,0.1906857109769646,0.01418120435697556,0,0,0.7848470606577962
5.202401853822494,1.495615476278033,0.1191748546277722,0.1619581450697751,0.3473887037929586,0.0007795445118418047,0.175491235802042,0.01192358401845794,0,0,0.7706965671103063
4.991035451832488,1.441267988981963,0.1089952125854552,0.2459753982704243,0.365796307107324,0.000780802229109395,0.1367689001963726,0.01081278822335374,0,0,0.8422055574621652
4.871728470037112,1.407347280808153,0.09563455489917061,0.4790218023010837,0.4181691898324047,0.0007811315644281613,0.08461431799794639,0.01026415122278657,0,0,1.011170216096644
4.83500880465499,1.40339474351282,0.08223365033970964,0.9131096264680565,0.4940337535870114,0.0007805737522911431,0.0420758044946794,0.01072812603060136,0,0,1.251779364725252
5.031627643970751,1.464437765624723,0.06744620694644254,1.891890550245648,0.6409981226065185,0.0007775358304290444,0.02500230937025849,0.01333659226125761,0,0,1.707873429777266
5.663603663768583,1.636139845533137,0.07300782904698622,2.224275313971731,0.7181090159483858,0.0007726865873697398,0.03281830044685519,0.01854060229383236,0,0,1.895197017569245
6.077296100170315,1.741715857843474,0.07705961331203781,2.390794747381719,0.7665380570823004,0.000769642109448282,0.03795447932891705,0.02207615676014552,0,0,2.010982655735347
5.856533976777541,1.678693843975964,0.0745603636229859,2.301860511110718,0.7437029678328884,0.0007712498272074986,0.03518443906473134,0.01990882957891671,0,0,1.959238721746175
5.022345505370311,1.452387582760707,0.0656626206233426,1.967538615152051,0.6548386594035378,0.0007774558731151467,0.02548343545092102,0.01297716825878936,0,0,1.754635322995938
4.125615602166818,1.210872665251486,0.05608168448596018,1.611287029704712,0.5640920550256467,0.0007843738772956623,0.01649207729639268,0.007186637996994705,0,0,1.549647283128024
"@

$lines = $csv -split "`n"
$rowIndex = 2
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split ","
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $colLetter = $columns[$i]
        $ws.Range("$colLetter$rowIndex").Value = $parts[$i]
    }
    $rowIndex++
}
